# "Generate Report for Handoff" — a new handoff has been kicked off for
# b.md, so the localization-status report needs to reflect that:
#   - Overview sheet: b.md's status flips from "Handed back: in sync with
#     en-US" to "Ready for handoff" for each locale column, and the
#     Latest Handoff Date advances to the new handoff timestamp.
#   - Per-locale sheets (zh-cn, de-de): b.md's row gets the new handoff
#     status, the new Latest Handoff File name/hyperlink, and the new
#     Latest Handoff Datetime.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$overviewDate = "2016-03-22 08:30:57"

$zhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate = "2016-03-22 08:30:53"

$deHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate = "2016-03-22 08:30:57"

# --- Overview sheet: row 3 is the b.md row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $overviewDate

# --- zh-cn sheet: row 3 is the b.md row ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = $zhHandoffFile
$wsZh.Range("E3").Value = $zhHandoffDate

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = $zhHandoffFile
    }
}

# --- de-de sheet: row 3 is the b.md row ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = $deHandoffFile
$wsDe.Range("E3").Value = $deHandoffDate

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = $deHandoffFile
    }
}
